$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing [Saucers]/4 block from column A to column B (rows 4-5),
# moving the cells (and their formatting) via cut/paste
$ws.Range("A4:A5").Cut($ws.Range("B4:B5"))

# Place the new [Time]/30 values into column A
$ws.Range("A4").Value2 = "[Time]"
$ws.Range("A5").Value2 = 30

# The cut left A4:A5 carrying the old formatting; reset it back to the
# workbook's default (unstyled) look
$ws.Range("A4:A5").Style = "Normal"

# Update the active selection to B6
$ws.Range("B6").Select()
